$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet contains a weekly price log for "Cebollín" at Feria Lagunitas de
# Puerto Montt. Two new weekly records need to be inserted into the existing
# chronological list (which is currently not perfectly sorted), causing all
# rows below each insertion point to shift down by one.
#
# 1) A new row is inserted right before the existing row 298 (fecha 45006).
# 2) A new row is inserted right before the existing row 350 (fecha 45005).
#
# Insert from the bottom up so the row numbers used for the second insertion
# are not invalidated by the first.

# --- Insert the second new row first (before original row 350) ---
$ws.Rows.Item(350).Insert()

$ws.Range("A350").Value = 4
$ws.Range("B350").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C350").Value = "Los Lagos"
$ws.Range("D350").Value = 45005
$ws.Range("E350").Value = 10
$ws.Range("F350").Value = 100112037
$ws.Range("G350").Value = "Cebollín"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 70
$ws.Range("K350").Value = 6500
$ws.Range("L350").Value = 7500
$ws.Range("M350").Value = 7000
$ws.Range("N350").Value = "`$/paquete 36 unidades"
$ws.Range("O350").Value = "Región Metropolitana"
$ws.Range("P350").Value = 194
$ws.Range("Q350").Value = 36
$ws.Range("R350").Value = "Hortaliza"

# --- Insert the first new row (before original row 298) ---
$ws.Rows.Item(298).Insert()

$ws.Range("A298").Value = 4
$ws.Range("B298").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C298").Value = "Los Lagos"
$ws.Range("D298").Value = 45006
$ws.Range("E298").Value = 10
$ws.Range("F298").Value = 100112037
$ws.Range("G298").Value = "Cebollín"
$ws.Range("H298").Value = "Sin especificar"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 180
$ws.Range("K298").Value = 6500
$ws.Range("L298").Value = 7000
$ws.Range("M298").Value = 6750
$ws.Range("N298").Value = "`$/paquete 36 unidades"
$ws.Range("O298").Value = "Región Metropolitana"
$ws.Range("P298").Value = 188
$ws.Range("Q298").Value = 36
$ws.Range("R298").Value = "Hortaliza"
